# Generate Report for Handback
# Applies the "handback" update to the localization-status workbook:
#  - Status text changes from "Ready for handoff" to "Handed back: in sync with en-US"
#  - zh-cn / de-de sheets get their "Latest Target File" (I) + "Latest Handback File" (J)
#    (and for de-de, "Latest Handback DateTime" (K)) columns populated, with a new
#    hyperlink on the Target File cell
#  - a couple of columns get wider to fit the new content

$wb = $excel.ActiveWorkbook

$targetMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4fb719931a18e171d93d74ae68f76a81bb3107cf/e2e/c7e37e7b-0198-4224-a154-a931b3cbc493.md"
$targetMdName = "c7e37e7b-0198-4224-a154-a931b3cbc493.md"

$statusText = "Handed back: in sync with en-US"
$zhXlf = "c7e37e7b-0198-4224-a154-a931b3cbc493.472148b32b812ac1221de1558083179dd5dd690a.zh-cn.xlf"
$deXlf = "c7e37e7b-0198-4224-a154-a931b3cbc493.472148b32b812ac1221de1558083179dd5dd690a.de-de.xlf"
$deHandbackDateTime = "2016-08-22 19:07:44"
$zhHandbackDateTime = "2016-08-22 19:07:37"

# ---- Overview sheet: widen zh-cn / de-de summary columns ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E1").ColumnWidth = 29.14
$wsOverview.Range("F1").ColumnWidth = 29.14

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C1").ColumnWidth = 29.14
$wsZh.Range("I1").ColumnWidth = 39.14
$wsZh.Range("J1").ColumnWidth = 39.14

$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

$wsZh.Range("K2").Value = $zhHandbackDateTime
$wsZh.Range("K3").Value = $zhHandbackDateTime

$wsZh.Range("J2").Value = $zhXlf
$wsZh.Range("J3").Value = $zhXlf

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $targetMdUrl, "", "", $targetMdName)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $targetMdUrl, "", "", $targetMdName)
$wsZh.Range("I2").Font.Underline = 2
$wsZh.Range("I2").Font.Color = 15570276
$wsZh.Range("I3").Font.Underline = 2
$wsZh.Range("I3").Font.Color = 15570276

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C1").ColumnWidth = 29.14
$wsDe.Range("I1").ColumnWidth = 39.14
$wsDe.Range("J1").ColumnWidth = 39.14

$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

$wsDe.Range("K2").Value = $deHandbackDateTime
$wsDe.Range("K3").Value = $deHandbackDateTime

$wsDe.Range("J2").Value = $deXlf
$wsDe.Range("J3").Value = $deXlf

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $targetMdUrl, "", "", $targetMdName)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $targetMdUrl, "", "", $targetMdName)
$wsDe.Range("I2").Font.Underline = 2
$wsDe.Range("I2").Font.Color = 15570276
$wsDe.Range("I3").Font.Underline = 2
$wsDe.Range("I3").Font.Color = 15570276
